# Auto-generated PowerShell Excel COM-interop script
# Applies the "Update countries & provincias Spain" edit to paises.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 20 de Marzo de 2020 a las 07:16"

# Full country table (rows 4-184), already updated with new case counts
# for Estados Unidos, Hungria, Guadalupe and Martinica, and re-sorted
# in descending order by "Casos totales" (column B).
$data = New-Object 'object[,]' 181,8
$data[0,0] = "China"
$data[0,1] = 80967
$data[0,2] = 39
$data[0,3] = 71150
$data[0,4] = 6569
$data[0,5] = 2136
$data[0,6] = 3
$data[0,7] = 3248
$data[1,0] = "Italia"
$data[1,1] = 41035
$data[1,2] = 0
$data[1,3] = 4440
$data[1,4] = 33190
$data[1,5] = 2498
$data[1,6] = 0
$data[1,7] = 3405
$data[2,0] = "Iran"
$data[2,1] = 18407
$data[2,2] = 0
$data[2,3] = 5979
$data[2,4] = 11144
$data[2,5] = 0
$data[2,6] = 0
$data[2,7] = 1284
$data[3,0] = "España"
$data[3,1] = 18077
$data[3,2] = 0
$data[3,3] = 1107
$data[3,4] = 16139
$data[3,5] = 939
$data[3,6] = 0
$data[3,7] = 831
$data[4,0] = "Alemania"
$data[4,1] = 15320
$data[4,2] = 0
$data[4,3] = 115
$data[4,4] = 15161
$data[4,5] = 2
$data[4,6] = 0
$data[4,7] = 44
$data[5,0] = "Estados Unidos"
$data[5,1] = 14366
$data[5,2] = 577
$data[5,3] = 125
$data[5,4] = 14024
$data[5,5] = 64
$data[5,6] = 10
$data[5,7] = 217
$data[6,0] = "Francia"
$data[6,1] = 10995
$data[6,2] = 0
$data[6,3] = 1295
$data[6,4] = 9328
$data[6,5] = 1122
$data[6,6] = 0
$data[6,7] = 372
$data[7,0] = "Corea del Sur"
$data[7,1] = 8652
$data[7,2] = 87
$data[7,3] = 2233
$data[7,4] = 6325
$data[7,5] = 59
$data[7,6] = 3
$data[7,7] = 94
$data[8,0] = "Suiza"
$data[8,1] = 4222
$data[8,2] = 0
$data[8,3] = 15
$data[8,4] = 4164
$data[8,5] = 0
$data[8,6] = 0
$data[8,7] = 43
$data[9,0] = "Reino Unido"
$data[9,1] = 3269
$data[9,2] = 0
$data[9,3] = 65
$data[9,4] = 3060
$data[9,5] = 20
$data[9,6] = 0
$data[9,7] = 144
$data[10,0] = "Paises Bajos"
$data[10,1] = 2460
$data[10,2] = 0
$data[10,3] = 2
$data[10,4] = 2382
$data[10,5] = 45
$data[10,6] = 0
$data[10,7] = 76
$data[11,0] = "Austria"
$data[11,1] = 2196
$data[11,2] = 17
$data[11,3] = 9
$data[11,4] = 2181
$data[11,5] = 13
$data[11,6] = 0
$data[11,7] = 6
$data[12,0] = "Belgica"
$data[12,1] = 1795
$data[12,2] = 0
$data[12,3] = 165
$data[12,4] = 1609
$data[12,5] = 130
$data[12,6] = 0
$data[12,7] = 21
$data[13,0] = "Noruega"
$data[13,1] = 1790
$data[13,2] = 0
$data[13,3] = 1
$data[13,4] = 1782
$data[13,5] = 27
$data[13,6] = 0
$data[13,7] = 7
$data[14,0] = "Suecia"
$data[14,1] = 1439
$data[14,2] = 0
$data[14,3] = 16
$data[14,4] = 1412
$data[14,5] = 21
$data[14,6] = 0
$data[14,7] = 11
$data[15,0] = "Dinamarca"
$data[15,1] = 1151
$data[15,2] = 0
$data[15,3] = 1
$data[15,4] = 1144
$data[15,5] = 30
$data[15,6] = 0
$data[15,7] = 6
$data[16,0] = "Japon"
$data[16,1] = 963
$data[16,2] = 20
$data[16,3] = 215
$data[16,4] = 715
$data[16,5] = 46
$data[16,6] = 0
$data[16,7] = 33
$data[17,0] = "Malasia"
$data[17,1] = 900
$data[17,2] = 0
$data[17,3] = 75
$data[17,4] = 823
$data[17,5] = 15
$data[17,6] = 0
$data[17,7] = 2
$data[18,0] = "Canada"
$data[18,1] = 873
$data[18,2] = 0
$data[18,3] = 11
$data[18,4] = 850
$data[18,5] = 1
$data[18,6] = 0
$data[18,7] = 12
$data[19,0] = "Australia"
$data[19,1] = 814
$data[19,2] = 58
$data[19,3] = 46
$data[19,4] = 761
$data[19,5] = 1
$data[19,6] = 0
$data[19,7] = 7
$data[20,0] = "Portugal"
$data[20,1] = 786
$data[20,2] = 0
$data[20,3] = 4
$data[20,4] = 778
$data[20,5] = 20
$data[20,6] = 0
$data[20,7] = 4
$data[21,0] = "Crucero"
$data[21,1] = 712
$data[21,2] = 0
$data[21,3] = 527
$data[21,4] = 178
$data[21,5] = 14
$data[21,6] = 0
$data[21,7] = 7
$data[22,0] = "Chequia"
$data[22,1] = 694
$data[22,2] = 0
$data[22,3] = 3
$data[22,4] = 691
$data[22,5] = 6
$data[22,6] = 0
$data[22,7] = 0
$data[23,0] = "Israel"
$data[23,1] = 677
$data[23,2] = 0
$data[23,3] = 14
$data[23,4] = 663
$data[23,5] = 6
$data[23,6] = 0
$data[23,7] = 0
$data[24,0] = "Brasil"
$data[24,1] = 647
$data[24,2] = 7
$data[24,3] = 2
$data[24,4] = 638
$data[24,5] = 18
$data[24,6] = 0
$data[24,7] = 7
$data[25,0] = "Irlanda"
$data[25,1] = 557
$data[25,2] = 0
$data[25,3] = 5
$data[25,4] = 549
$data[25,5] = 6
$data[25,6] = 0
$data[25,7] = 3
$data[26,0] = "Grecia"
$data[26,1] = 464
$data[26,2] = 0
$data[26,3] = 19
$data[26,4] = 439
$data[26,5] = 16
$data[26,6] = 0
$data[26,7] = 6
$data[27,0] = "Catar"
$data[27,1] = 460
$data[27,2] = 0
$data[27,3] = 10
$data[27,4] = 450
$data[27,5] = 6
$data[27,6] = 0
$data[27,7] = 0
$data[28,0] = "Pakistan"
$data[28,1] = 454
$data[28,2] = 0
$data[28,3] = 13
$data[28,4] = 439
$data[28,5] = 0
$data[28,6] = 0
$data[28,7] = 2
$data[29,0] = "Finlandia"
$data[29,1] = 400
$data[29,2] = 0
$data[29,3] = 10
$data[29,4] = 390
$data[29,5] = 2
$data[29,6] = 0
$data[29,7] = 0
$data[30,0] = "Turquia"
$data[30,1] = 359
$data[30,2] = 0
$data[30,3] = 0
$data[30,4] = 355
$data[30,5] = 0
$data[30,6] = 0
$data[30,7] = 4
$data[31,0] = "Polonia"
$data[31,1] = 355
$data[31,2] = 0
$data[31,3] = 13
$data[31,4] = 337
$data[31,5] = 3
$data[31,6] = 0
$data[31,7] = 5
$data[32,0] = "Singapur"
$data[32,1] = 345
$data[32,2] = 0
$data[32,3] = 124
$data[32,4] = 221
$data[32,5] = 14
$data[32,6] = 0
$data[32,7] = 0
$data[33,0] = "Chile"
$data[33,1] = 342
$data[33,2] = 0
$data[33,3] = 0
$data[33,4] = 342
$data[33,5] = 6
$data[33,6] = 0
$data[33,7] = 0
$data[34,0] = "Luxemburgo"
$data[34,1] = 335
$data[34,2] = 0
$data[34,3] = 6
$data[34,4] = 325
$data[34,5] = 1
$data[34,6] = 0
$data[34,7] = 4
$data[35,0] = "Islandia"
$data[35,1] = 330
$data[35,2] = 0
$data[35,3] = 5
$data[35,4] = 325
$data[35,5] = 1
$data[35,6] = 0
$data[35,7] = 0
$data[36,0] = "Tailandia"
$data[36,1] = 322
$data[36,2] = 50
$data[36,3] = 42
$data[36,4] = 279
$data[36,5] = 1
$data[36,6] = 0
$data[36,7] = 1
$data[37,0] = "Eslovenia"
$data[37,1] = 319
$data[37,2] = 0
$data[37,3] = 0
$data[37,4] = 318
$data[37,5] = 6
$data[37,6] = 0
$data[37,7] = 1
$data[38,0] = "Indonesia"
$data[38,1] = 309
$data[38,2] = 0
$data[38,3] = 15
$data[38,4] = 269
$data[38,5] = 0
$data[38,6] = 0
$data[38,7] = 25
$data[39,0] = "Barein"
$data[39,1] = 279
$data[39,2] = 0
$data[39,3] = 110
$data[39,4] = 168
$data[39,5] = 4
$data[39,6] = 0
$data[39,7] = 1
$data[40,0] = "Rumania"
$data[40,1] = 277
$data[40,2] = 0
$data[40,3] = 25
$data[40,4] = 252
$data[40,5] = 5
$data[40,6] = 0
$data[40,7] = 0
$data[41,0] = "Arabia Saudita"
$data[41,1] = 274
$data[41,2] = 0
$data[41,3] = 8
$data[41,4] = 266
$data[41,5] = 0
$data[41,6] = 0
$data[41,7] = 0
$data[42,0] = "Estonia"
$data[42,1] = 267
$data[42,2] = 0
$data[42,3] = 1
$data[42,4] = 266
$data[42,5] = 1
$data[42,6] = 0
$data[42,7] = 0
$data[43,0] = "Ecuador"
$data[43,1] = 260
$data[43,2] = 0
$data[43,3] = 1
$data[43,4] = 256
$data[43,5] = 2
$data[43,6] = 0
$data[43,7] = 3
$data[44,0] = "Egipto"
$data[44,1] = 256
$data[44,2] = 0
$data[44,3] = 42
$data[44,4] = 207
$data[44,5] = 0
$data[44,6] = 0
$data[44,7] = 7
$data[45,0] = "Peru"
$data[45,1] = 234
$data[45,2] = 0
$data[45,3] = 1
$data[45,4] = 230
$data[45,5] = 7
$data[45,6] = 2
$data[45,7] = 3
$data[46,0] = "Filipinas"
$data[46,1] = 217
$data[46,2] = 0
$data[46,3] = 8
$data[46,4] = 192
$data[46,5] = 1
$data[46,6] = 0
$data[46,7] = 17
$data[47,0] = "Hong Kong"
$data[47,1] = 208
$data[47,2] = 0
$data[47,3] = 98
$data[47,4] = 106
$data[47,5] = 4
$data[47,6] = 0
$data[47,7] = 4
$data[48,0] = "India"
$data[48,1] = 201
$data[48,2] = 7
$data[48,3] = 20
$data[48,4] = 176
$data[48,5] = 0
$data[48,6] = 1
$data[48,7] = 5
$data[49,0] = "Rusia"
$data[49,1] = 199
$data[49,2] = 0
$data[49,3] = 8
$data[49,4] = 190
$data[49,5] = 0
$data[49,6] = 0
$data[49,7] = 1
$data[50,0] = "Irak"
$data[50,1] = 192
$data[50,2] = 0
$data[50,3] = 49
$data[50,4] = 130
$data[50,5] = 0
$data[50,6] = 0
$data[50,7] = 13
$data[51,0] = "Mexico"
$data[51,1] = 164
$data[51,2] = 46
$data[51,3] = 4
$data[51,4] = 159
$data[51,5] = 1
$data[51,6] = 0
$data[51,7] = 1
$data[52,0] = "Libano"
$data[52,1] = 157
$data[52,2] = 0
$data[52,3] = 4
$data[52,4] = 149
$data[52,5] = 3
$data[52,6] = 0
$data[52,7] = 4
$data[53,0] = "Sudafrica"
$data[53,1] = 150
$data[53,2] = 0
$data[53,3] = 0
$data[53,4] = 150
$data[53,5] = 0
$data[53,6] = 0
$data[53,7] = 0
$data[54,0] = "Kuwait"
$data[54,1] = 148
$data[54,2] = 0
$data[54,3] = 18
$data[54,4] = 130
$data[54,5] = 5
$data[54,6] = 0
$data[54,7] = 0
$data[55,0] = "San Marino"
$data[55,1] = 144
$data[55,2] = 0
$data[55,3] = 4
$data[55,4] = 126
$data[55,5] = 12
$data[55,6] = 0
$data[55,7] = 14
$data[56,0] = "Emiratos Arabes Unidos"
$data[56,1] = 140
$data[56,2] = 0
$data[56,3] = 31
$data[56,4] = 109
$data[56,5] = 2
$data[56,6] = 0
$data[56,7] = 0
$data[57,0] = "Panama"
$data[57,1] = 137
$data[57,2] = 0
$data[57,3] = 1
$data[57,4] = 135
$data[57,5] = 7
$data[57,6] = 0
$data[57,7] = 1
$data[58,0] = "Colombia"
$data[58,1] = 128
$data[58,2] = 20
$data[58,3] = 1
$data[58,4] = 127
$data[58,5] = 0
$data[58,6] = 0
$data[58,7] = 0
$data[59,0] = "Argentina"
$data[59,1] = 128
$data[59,2] = 0
$data[59,3] = 3
$data[59,4] = 122
$data[59,5] = 0
$data[59,6] = 0
$data[59,7] = 3
$data[60,0] = "Eslovaquia"
$data[60,1] = 124
$data[60,2] = 0
$data[60,3] = 0
$data[60,4] = 124
$data[60,5] = 2
$data[60,6] = 0
$data[60,7] = 0
$data[61,0] = "Armenia"
$data[61,1] = 122
$data[61,2] = 0
$data[61,3] = 1
$data[61,4] = 121
$data[61,5] = 2
$data[61,6] = 0
$data[61,7] = 0
$data[62,0] = "Croacia"
$data[62,1] = 110
$data[62,2] = 0
$data[62,3] = 5
$data[62,4] = 104
$data[62,5] = 0
$data[62,6] = 0
$data[62,7] = 1
$data[63,0] = "Taiwan"
$data[63,1] = 108
$data[63,2] = 0
$data[63,3] = 26
$data[63,4] = 81
$data[63,5] = 0
$data[63,6] = 0
$data[63,7] = 1
$data[64,0] = "Bulgaria"
$data[64,1] = 107
$data[64,2] = 0
$data[64,3] = 0
$data[64,4] = 104
$data[64,5] = 0
$data[64,6] = 0
$data[64,7] = 3
$data[65,0] = "Serbia"
$data[65,1] = 103
$data[65,2] = 0
$data[65,3] = 1
$data[65,4] = 102
$data[65,5] = 4
$data[65,6] = 0
$data[65,7] = 0
$data[66,0] = "Uruguay"
$data[66,1] = 94
$data[66,2] = 15
$data[66,3] = 0
$data[66,4] = 94
$data[66,5] = 0
$data[66,6] = 0
$data[66,7] = 0
$data[67,0] = "Argelia"
$data[67,1] = 90
$data[67,2] = 0
$data[67,3] = 32
$data[67,4] = 49
$data[67,5] = 0
$data[67,6] = 0
$data[67,7] = 9
$data[68,0] = "Costa Rica"
$data[68,1] = 89
$data[68,2] = 2
$data[68,3] = 0
$data[68,4] = 87
$data[68,5] = 2
$data[68,6] = 1
$data[68,7] = 2
$data[69,0] = "Letonia"
$data[69,1] = 86
$data[69,2] = 0
$data[69,3] = 1
$data[69,4] = 85
$data[69,5] = 0
$data[69,6] = 0
$data[69,7] = 0
$data[70,0] = "Hungria"
$data[70,1] = 85
$data[70,2] = 12
$data[70,3] = 2
$data[70,4] = 82
$data[70,5] = 4
$data[70,6] = 0
$data[70,7] = 1
$data[71,0] = "Vietnam"
$data[71,1] = 85
$data[71,2] = 0
$data[71,3] = 16
$data[71,4] = 69
$data[71,5] = 0
$data[71,6] = 0
$data[71,7] = 0
$data[72,0] = "Principado de Andorra"
$data[72,1] = 74
$data[72,2] = 0
$data[72,3] = 1
$data[72,4] = 73
$data[72,5] = 0
$data[72,6] = 0
$data[72,7] = 0
$data[73,0] = "Brunei"
$data[73,1] = 73
$data[73,2] = 0
$data[73,3] = 0
$data[73,4] = 73
$data[73,5] = 2
$data[73,6] = 0
$data[73,7] = 0
$data[74,0] = "Islas Feroe"
$data[74,1] = 72
$data[74,2] = 0
$data[74,3] = 1
$data[74,4] = 71
$data[74,5] = 0
$data[74,6] = 0
$data[74,7] = 0
$data[75,0] = "Jordania"
$data[75,1] = 69
$data[75,2] = 0
$data[75,3] = 1
$data[75,4] = 68
$data[75,5] = 0
$data[75,6] = 0
$data[75,7] = 0
$data[76,0] = "Republica de Chipre"
$data[76,1] = 67
$data[76,2] = 0
$data[76,3] = 0
$data[76,4] = 67
$data[76,5] = 1
$data[76,6] = 0
$data[76,7] = 0
$data[77,0] = "Albania"
$data[77,1] = 64
$data[77,2] = 0
$data[77,3] = 0
$data[77,4] = 62
$data[77,5] = 2
$data[77,6] = 0
$data[77,7] = 2
$data[78,0] = "Bosnia y Herzegovina"
$data[78,1] = 64
$data[78,2] = 0
$data[78,3] = 2
$data[78,4] = 62
$data[78,5] = 0
$data[78,6] = 0
$data[78,7] = 0
$data[79,0] = "Marruecos"
$data[79,1] = 63
$data[79,2] = 0
$data[79,3] = 2
$data[79,4] = 59
$data[79,5] = 1
$data[79,6] = 0
$data[79,7] = 2
$data[80,0] = "Sri Lanka"
$data[80,1] = 60
$data[80,2] = 0
$data[80,3] = 3
$data[80,4] = 57
$data[80,5] = 0
$data[80,6] = 0
$data[80,7] = 0
$data[81,0] = "Malta"
$data[81,1] = 53
$data[81,2] = 0
$data[81,3] = 2
$data[81,4] = 51
$data[81,5] = 0
$data[81,6] = 0
$data[81,7] = 0
$data[82,0] = "Bielorrusia"
$data[82,1] = 51
$data[82,2] = 0
$data[82,3] = 5
$data[82,4] = 46
$data[82,5] = 0
$data[82,6] = 0
$data[82,7] = 0
$data[83,0] = "Republica de Macedonia"
$data[83,1] = 50
$data[83,2] = 0
$data[83,3] = 1
$data[83,4] = 49
$data[83,5] = 1
$data[83,6] = 0
$data[83,7] = 0
$data[84,0] = "Kazajistan"
$data[84,1] = 49
$data[84,2] = 5
$data[84,3] = 0
$data[84,4] = 49
$data[84,5] = 0
$data[84,6] = 0
$data[84,7] = 0
$data[85,0] = "Moldavia"
$data[85,1] = 49
$data[85,2] = 0
$data[85,3] = 1
$data[85,4] = 47
$data[85,5] = 3
$data[85,6] = 0
$data[85,7] = 1
$data[86,0] = "Lituania"
$data[86,1] = 48
$data[86,2] = 0
$data[86,3] = 1
$data[86,4] = 47
$data[86,5] = 1
$data[86,6] = 0
$data[86,7] = 0
$data[87,0] = "Oman"
$data[87,1] = 48
$data[87,2] = 0
$data[87,3] = 13
$data[87,4] = 35
$data[87,5] = 0
$data[87,6] = 0
$data[87,7] = 0
$data[88,0] = "Estado de Palestina"
$data[88,1] = 47
$data[88,2] = 0
$data[88,3] = 0
$data[88,4] = 47
$data[88,5] = 0
$data[88,6] = 0
$data[88,7] = 0
$data[89,0] = "Guadalupe"
$data[89,1] = 45
$data[89,2] = 12
$data[89,3] = 0
$data[89,4] = 45
$data[89,5] = 0
$data[89,6] = 0
$data[89,7] = 0
$data[90,0] = "Azerbaiyan"
$data[90,1] = 44
$data[90,2] = 0
$data[90,3] = 7
$data[90,4] = 36
$data[90,5] = 0
$data[90,6] = 0
$data[90,7] = 1
$data[91,0] = "Venezuela"
$data[91,1] = 42
$data[91,2] = 0
$data[91,3] = 0
$data[91,4] = 42
$data[91,5] = 0
$data[91,6] = 0
$data[91,7] = 0
$data[92,0] = "Georgia"
$data[92,1] = 40
$data[92,2] = 0
$data[92,3] = 1
$data[92,4] = 39
$data[92,5] = 1
$data[92,6] = 0
$data[92,7] = 0
$data[93,0] = "Nueva Zelanda"
$data[93,1] = 39
$data[93,2] = 11
$data[93,3] = 0
$data[93,4] = 39
$data[93,5] = 0
$data[93,6] = 0
$data[93,7] = 0
$data[94,0] = "Tunez"
$data[94,1] = 39
$data[94,2] = 0
$data[94,3] = 1
$data[94,4] = 37
$data[94,5] = 2
$data[94,6] = 0
$data[94,7] = 1
$data[95,0] = "Camboya"
$data[95,1] = 37
$data[95,2] = 0
$data[95,3] = 1
$data[95,4] = 36
$data[95,5] = 0
$data[95,6] = 0
$data[95,7] = 0
$data[96,0] = "Senegal"
$data[96,1] = 36
$data[96,2] = 0
$data[96,3] = 2
$data[96,4] = 34
$data[96,5] = 0
$data[96,6] = 0
$data[96,7] = 0
$data[97,0] = "Republica Dominicana"
$data[97,1] = 34
$data[97,2] = 0
$data[97,3] = 0
$data[97,4] = 32
$data[97,5] = 0
$data[97,6] = 0
$data[97,7] = 2
$data[98,0] = "Burkina Faso"
$data[98,1] = 33
$data[98,2] = 0
$data[98,3] = 0
$data[98,4] = 32
$data[98,5] = 0
$data[98,6] = 0
$data[98,7] = 1
$data[99,0] = "Reunion"
$data[99,1] = 28
$data[99,2] = 0
$data[99,3] = 0
$data[99,4] = 28
$data[99,5] = 0
$data[99,6] = 0
$data[99,7] = 0
$data[100,0] = "Liechtenstein"
$data[100,1] = 28
$data[100,2] = 0
$data[100,3] = 0
$data[100,4] = 28
$data[100,5] = 0
$data[100,6] = 0
$data[100,7] = 0
$data[101,0] = "Ucrania"
$data[101,1] = 26
$data[101,2] = 0
$data[101,3] = 0
$data[101,4] = 23
$data[101,5] = 0
$data[101,6] = 0
$data[101,7] = 3
$data[102,0] = "Honduras"
$data[102,1] = 24
$data[102,2] = 12
$data[102,3] = 0
$data[102,4] = 24
$data[102,5] = 0
$data[102,6] = 0
$data[102,7] = 0
$data[103,0] = "Uzbekistan"
$data[103,1] = 23
$data[103,2] = 0
$data[103,3] = 0
$data[103,4] = 23
$data[103,5] = 0
$data[103,6] = 0
$data[103,7] = 0
$data[104,0] = "Martinica"
$data[104,1] = 23
$data[104,2] = 0
$data[104,3] = 0
$data[104,4] = 22
$data[104,5] = 7
$data[104,6] = 0
$data[104,7] = 1
$data[105,0] = "Afganistan"
$data[105,1] = 22
$data[105,2] = 0
$data[105,3] = 1
$data[105,4] = 21
$data[105,5] = 0
$data[105,6] = 0
$data[105,7] = 0
$data[106,0] = "Banglades"
$data[106,1] = 18
$data[106,2] = 0
$data[106,3] = 3
$data[106,4] = 14
$data[106,5] = 0
$data[106,6] = 0
$data[106,7] = 1
$data[107,0] = "Bolivia"
$data[107,1] = 17
$data[107,2] = 2
$data[107,3] = 0
$data[107,4] = 17
$data[107,5] = 0
$data[107,6] = 0
$data[107,7] = 0
$data[108,0] = "Macao"
$data[108,1] = 17
$data[108,2] = 0
$data[108,3] = 10
$data[108,4] = 7
$data[108,5] = 0
$data[108,6] = 0
$data[108,7] = 0
$data[109,0] = "Cuba"
$data[109,1] = 16
$data[109,2] = 5
$data[109,3] = 0
$data[109,4] = 15
$data[109,5] = 0
$data[109,6] = 0
$data[109,7] = 1
$data[110,0] = "Jamaica"
$data[110,1] = 16
$data[110,2] = 1
$data[110,3] = 2
$data[110,4] = 13
$data[110,5] = 0
$data[110,6] = 0
$data[110,7] = 1
$data[111,0] = "Guayana Francesa"
$data[111,1] = 15
$data[111,2] = 0
$data[111,3] = 0
$data[111,4] = 15
$data[111,5] = 0
$data[111,6] = 0
$data[111,7] = 0
$data[112,0] = "Consejo Danes para los Refugiados"
$data[112,1] = 14
$data[112,2] = 0
$data[112,3] = 0
$data[112,4] = 14
$data[112,5] = 0
$data[112,6] = 0
$data[112,7] = 0
$data[113,0] = "Camerun"
$data[113,1] = 13
$data[113,2] = 0
$data[113,3] = 0
$data[113,4] = 13
$data[113,5] = 0
$data[113,6] = 0
$data[113,7] = 0
$data[114,0] = "Maldivas"
$data[114,1] = 13
$data[114,2] = 0
$data[114,3] = 0
$data[114,4] = 13
$data[114,5] = 0
$data[114,6] = 0
$data[114,7] = 0
$data[115,0] = "Montenegro"
$data[115,1] = 13
$data[115,2] = 0
$data[115,3] = 0
$data[115,4] = 13
$data[115,5] = 0
$data[115,6] = 0
$data[115,7] = 0
$data[116,0] = "Paraguay"
$data[116,1] = 13
$data[116,2] = 0
$data[116,3] = 0
$data[116,4] = 13
$data[116,5] = 1
$data[116,6] = 0
$data[116,7] = 0
$data[117,0] = "Guam"
$data[117,1] = 12
$data[117,2] = 0
$data[117,3] = 0
$data[117,4] = 12
$data[117,5] = 0
$data[117,6] = 0
$data[117,7] = 0
$data[118,0] = "Nigeria"
$data[118,1] = 12
$data[118,2] = 0
$data[118,3] = 1
$data[118,4] = 11
$data[118,5] = 0
$data[118,6] = 0
$data[118,7] = 0
$data[119,0] = "Polinesia Francesa"
$data[119,1] = 11
$data[119,2] = 5
$data[119,3] = 0
$data[119,4] = 11
$data[119,5] = 0
$data[119,6] = 0
$data[119,7] = 0
$data[120,0] = "Ruanda"
$data[120,1] = 11
$data[120,2] = 0
$data[120,3] = 0
$data[120,4] = 11
$data[120,5] = 0
$data[120,6] = 0
$data[120,7] = 0
$data[121,0] = "Ghana"
$data[121,1] = 11
$data[121,2] = 0
$data[121,3] = 0
$data[121,4] = 11
$data[121,5] = 0
$data[121,6] = 0
$data[121,7] = 0
$data[122,0] = "Monaco"
$data[122,1] = 10
$data[122,2] = 0
$data[122,3] = 0
$data[122,4] = 10
$data[122,5] = 0
$data[122,6] = 0
$data[122,7] = 0
$data[123,0] = "Gibraltar"
$data[123,1] = 10
$data[123,2] = 0
$data[123,3] = 2
$data[123,4] = 8
$data[123,5] = 0
$data[123,6] = 0
$data[123,7] = 0
$data[124,0] = "Etiopia"
$data[124,1] = 9
$data[124,2] = 2
$data[124,3] = 0
$data[124,4] = 9
$data[124,5] = 0
$data[124,6] = 0
$data[124,7] = 0
$data[125,0] = "Trinidad yTobago"
$data[125,1] = 9
$data[125,2] = 0
$data[125,3] = 0
$data[125,4] = 9
$data[125,5] = 0
$data[125,6] = 0
$data[125,7] = 0
$data[126,0] = "Guatemala"
$data[126,1] = 9
$data[126,2] = 0
$data[126,3] = 0
$data[126,4] = 8
$data[126,5] = 0
$data[126,6] = 0
$data[126,7] = 1
$data[127,0] = "Costa de Marfil"
$data[127,1] = 9
$data[127,2] = 0
$data[127,3] = 1
$data[127,4] = 8
$data[127,5] = 0
$data[127,6] = 0
$data[127,7] = 0
$data[128,0] = "Mauricio"
$data[128,1] = 7
$data[128,2] = 0
$data[128,3] = 0
$data[128,4] = 7
$data[128,5] = 0
$data[128,6] = 0
$data[128,7] = 0
$data[129,0] = "Kenia"
$data[129,1] = 7
$data[129,2] = 0
$data[129,3] = 0
$data[129,4] = 7
$data[129,5] = 0
$data[129,6] = 0
$data[129,7] = 0
$data[130,0] = "Kirguistan"
$data[130,1] = 6
$data[130,2] = 3
$data[130,3] = 0
$data[130,4] = 6
$data[130,5] = 0
$data[130,6] = 0
$data[130,7] = 0
$data[131,0] = "Guinea Ecuatorial"
$data[131,1] = 6
$data[131,2] = 0
$data[131,3] = 0
$data[131,4] = 6
$data[131,5] = 0
$data[131,6] = 0
$data[131,7] = 0
$data[132,0] = "Seychelles"
$data[132,1] = 6
$data[132,2] = 0
$data[132,3] = 0
$data[132,4] = 6
$data[132,5] = 0
$data[132,6] = 0
$data[132,7] = 0
$data[133,0] = "Puerto Rico"
$data[133,1] = 6
$data[133,2] = 0
$data[133,3] = 0
$data[133,4] = 6
$data[133,5] = 0
$data[133,6] = 0
$data[133,7] = 0
$data[134,0] = "Mongolia"
$data[134,1] = 6
$data[134,2] = 0
$data[134,3] = 0
$data[134,4] = 6
$data[134,5] = 0
$data[134,6] = 0
$data[134,7] = 0
$data[135,0] = "Tanzania"
$data[135,1] = 6
$data[135,2] = 0
$data[135,3] = 0
$data[135,4] = 6
$data[135,5] = 0
$data[135,6] = 0
$data[135,7] = 0
$data[136,0] = "Barbados"
$data[136,1] = 5
$data[136,2] = 0
$data[136,3] = 0
$data[136,4] = 5
$data[136,5] = 0
$data[136,6] = 0
$data[136,7] = 0
$data[137,0] = "Guyana"
$data[137,1] = 5
$data[137,2] = 0
$data[137,3] = 0
$data[137,4] = 4
$data[137,5] = 0
$data[137,6] = 0
$data[137,7] = 1
$data[138,0] = "Aruba"
$data[138,1] = 5
$data[138,2] = 0
$data[138,3] = 1
$data[138,4] = 4
$data[138,5] = 0
$data[138,6] = 0
$data[138,7] = 0
$data[139,0] = "Mayotte"
$data[139,1] = 4
$data[139,2] = 0
$data[139,3] = 0
$data[139,4] = 4
$data[139,5] = 0
$data[139,6] = 0
$data[139,7] = 0
$data[140,0] = "San Bartolome"
$data[140,1] = 3
$data[140,2] = 0
$data[140,3] = 0
$data[140,4] = 3
$data[140,5] = 0
$data[140,6] = 0
$data[140,7] = 0
$data[141,0] = "Gabon"
$data[141,1] = 3
$data[141,2] = 0
$data[141,3] = 0
$data[141,4] = 3
$data[141,5] = 0
$data[141,6] = 0
$data[141,7] = 0
$data[142,0] = "Congo"
$data[142,1] = 3
$data[142,2] = 0
$data[142,3] = 0
$data[142,4] = 3
$data[142,5] = 0
$data[142,6] = 0
$data[142,7] = 0
$data[143,0] = "San Martin (Parte Francesa)"
$data[143,1] = 3
$data[143,2] = 0
$data[143,3] = 0
$data[143,4] = 3
$data[143,5] = 0
$data[143,6] = 0
$data[143,7] = 0
$data[144,0] = "Islas Virgenes de los Estados Unidos"
$data[144,1] = 3
$data[144,2] = 0
$data[144,3] = 0
$data[144,4] = 3
$data[144,5] = 0
$data[144,6] = 0
$data[144,7] = 0
$data[145,0] = "Bahamas"
$data[145,1] = 3
$data[145,2] = 0
$data[145,3] = 0
$data[145,4] = 3
$data[145,5] = 0
$data[145,6] = 0
$data[145,7] = 0
$data[146,0] = "Namibia"
$data[146,1] = 3
$data[146,2] = 0
$data[146,3] = 0
$data[146,4] = 3
$data[146,5] = 0
$data[146,6] = 0
$data[146,7] = 0
$data[147,0] = "Islas Caimanes"
$data[147,1] = 3
$data[147,2] = 0
$data[147,3] = 0
$data[147,4] = 2
$data[147,5] = 0
$data[147,6] = 0
$data[147,7] = 1
$data[148,0] = "Curazao"
$data[148,1] = 3
$data[148,2] = 0
$data[148,3] = 0
$data[148,4] = 2
$data[148,5] = 0
$data[148,6] = 0
$data[148,7] = 1
$data[149,0] = "Groenlandia"
$data[149,1] = 2
$data[149,2] = 0
$data[149,3] = 0
$data[149,4] = 2
$data[149,5] = 0
$data[149,6] = 0
$data[149,7] = 0
$data[150,0] = "Santa Lucia"
$data[150,1] = 2
$data[150,2] = 0
$data[150,3] = 0
$data[150,4] = 2
$data[150,5] = 0
$data[150,6] = 0
$data[150,7] = 0
$data[151,0] = "Benin"
$data[151,1] = 2
$data[151,2] = 0
$data[151,3] = 0
$data[151,4] = 2
$data[151,5] = 0
$data[151,6] = 0
$data[151,7] = 0
$data[152,0] = "Mauritania"
$data[152,1] = 2
$data[152,2] = 0
$data[152,3] = 0
$data[152,4] = 2
$data[152,5] = 0
$data[152,6] = 0
$data[152,7] = 0
$data[153,0] = "Zambia"
$data[153,1] = 2
$data[153,2] = 0
$data[153,3] = 0
$data[153,4] = 2
$data[153,5] = 0
$data[153,6] = 0
$data[153,7] = 0
$data[154,0] = "Butan"
$data[154,1] = 2
$data[154,2] = 1
$data[154,3] = 0
$data[154,4] = 2
$data[154,5] = 0
$data[154,6] = 0
$data[154,7] = 0
$data[155,0] = "Nueva Caledonia"
$data[155,1] = 2
$data[155,2] = 0
$data[155,3] = 0
$data[155,4] = 2
$data[155,5] = 0
$data[155,6] = 0
$data[155,7] = 0
$data[156,0] = "Bermudas"
$data[156,1] = 2
$data[156,2] = 0
$data[156,3] = 0
$data[156,4] = 2
$data[156,5] = 0
$data[156,6] = 0
$data[156,7] = 0
$data[157,0] = "Liberia"
$data[157,1] = 2
$data[157,2] = 0
$data[157,3] = 0
$data[157,4] = 2
$data[157,5] = 0
$data[157,6] = 0
$data[157,7] = 0
$data[158,0] = "Haiti"
$data[158,1] = 2
$data[158,2] = 2
$data[158,3] = 0
$data[158,4] = 2
$data[158,5] = 0
$data[158,6] = 0
$data[158,7] = 0
$data[159,0] = "Sudan"
$data[159,1] = 2
$data[159,2] = 0
$data[159,3] = 0
$data[159,4] = 1
$data[159,5] = 0
$data[159,6] = 0
$data[159,7] = 1
$data[160,0] = "Republica de Yibuti"
$data[160,1] = 1
$data[160,2] = 0
$data[160,3] = 0
$data[160,4] = 1
$data[160,5] = 0
$data[160,6] = 0
$data[160,7] = 0
$data[161,0] = "Republica de Africa Central"
$data[161,1] = 1
$data[161,2] = 0
$data[161,3] = 0
$data[161,4] = 1
$data[161,5] = 0
$data[161,6] = 0
$data[161,7] = 0
$data[162,0] = "San Martin (Parte Holandesa)"
$data[162,1] = 1
$data[162,2] = 0
$data[162,3] = 0
$data[162,4] = 1
$data[162,5] = 0
$data[162,6] = 0
$data[162,7] = 0
$data[163,0] = "Fiyi"
$data[163,1] = 1
$data[163,2] = 0
$data[163,3] = 0
$data[163,4] = 1
$data[163,5] = 0
$data[163,6] = 0
$data[163,7] = 0
$data[164,0] = "El Salvador"
$data[164,1] = 1
$data[164,2] = 0
$data[164,3] = 0
$data[164,4] = 1
$data[164,5] = 0
$data[164,6] = 0
$data[164,7] = 0
$data[165,0] = "Nicaragua"
$data[165,1] = 1
$data[165,2] = 0
$data[165,3] = 0
$data[165,4] = 1
$data[165,5] = 0
$data[165,6] = 0
$data[165,7] = 0
$data[166,0] = "Somalia"
$data[166,1] = 1
$data[166,2] = 0
$data[166,3] = 0
$data[166,4] = 1
$data[166,5] = 0
$data[166,6] = 0
$data[166,7] = 0
$data[167,0] = "Surinam"
$data[167,1] = 1
$data[167,2] = 0
$data[167,3] = 0
$data[167,4] = 1
$data[167,5] = 0
$data[167,6] = 0
$data[167,7] = 0
$data[168,0] = "Cabo Verde"
$data[168,1] = 1
$data[168,2] = 1
$data[168,3] = 0
$data[168,4] = 1
$data[168,5] = 0
$data[168,6] = 0
$data[168,7] = 0
$data[169,0] = "Togo"
$data[169,1] = 1
$data[169,2] = 0
$data[169,3] = 0
$data[169,4] = 1
$data[169,5] = 0
$data[169,6] = 0
$data[169,7] = 0
$data[170,0] = "Montserrat"
$data[170,1] = 1
$data[170,2] = 0
$data[170,3] = 0
$data[170,4] = 1
$data[170,5] = 0
$data[170,6] = 0
$data[170,7] = 0
$data[171,0] = "Gambia"
$data[171,1] = 1
$data[171,2] = 0
$data[171,3] = 0
$data[171,4] = 1
$data[171,5] = 0
$data[171,6] = 0
$data[171,7] = 0
$data[172,0] = "San Vicente y las Granadinas"
$data[172,1] = 1
$data[172,2] = 0
$data[172,3] = 0
$data[172,4] = 1
$data[172,5] = 0
$data[172,6] = 0
$data[172,7] = 0
$data[173,0] = "Republica del Chad"
$data[173,1] = 1
$data[173,2] = 0
$data[173,3] = 0
$data[173,4] = 1
$data[173,5] = 0
$data[173,6] = 0
$data[173,7] = 0
$data[174,0] = "Suazilandia"
$data[174,1] = 1
$data[174,2] = 0
$data[174,3] = 0
$data[174,4] = 1
$data[174,5] = 0
$data[174,6] = 0
$data[174,7] = 0
$data[175,0] = "Niger"
$data[175,1] = 1
$data[175,2] = 0
$data[175,3] = 0
$data[175,4] = 1
$data[175,5] = 0
$data[175,6] = 0
$data[175,7] = 0
$data[176,0] = "Antigua y Barbuda"
$data[176,1] = 1
$data[176,2] = 0
$data[176,3] = 0
$data[176,4] = 1
$data[176,5] = 0
$data[176,6] = 0
$data[176,7] = 0
$data[177,0] = "Santa Sede"
$data[177,1] = 1
$data[177,2] = 0
$data[177,3] = 0
$data[177,4] = 1
$data[177,5] = 0
$data[177,6] = 0
$data[177,7] = 0
$data[178,0] = "Isla de Man"
$data[178,1] = 1
$data[178,2] = 0
$data[178,3] = 0
$data[178,4] = 1
$data[178,5] = 0
$data[178,6] = 0
$data[178,7] = 0
$data[179,0] = "Guinea"
$data[179,1] = 1
$data[179,2] = 0
$data[179,3] = 0
$data[179,4] = 1
$data[179,5] = 0
$data[179,6] = 0
$data[179,7] = 0
$data[180,0] = "Nepal"
$data[180,1] = 1
$data[180,2] = 0
$data[180,3] = 1
$data[180,4] = 0
$data[180,5] = 0
$data[180,6] = 0
$data[180,7] = 0

$ws.Range("A4:H184").Value = $data

